$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.568.70'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -0.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.848.20'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  -0.39%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '262.53'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -1.03%  '

$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5337'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +1.91%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3144'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  -4.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06908'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +1.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.75'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -0.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7656'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  -1.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07824'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +1.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.853.67'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -0.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '89.60'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +0.98%  '

$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("E17").Value = '  +0.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007953'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.619.71'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -0.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.085.77'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.638'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.014'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.342'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -2.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.206'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -0.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.60'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -1.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.688'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +0.37%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.61'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -0.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.298'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +2.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08784'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +0.26%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.106'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -1.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04852'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +0.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7356'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +2.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.931'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +2.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.137'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -0.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.112'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -0.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.332'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +5.72%  '

$ws.Range("E39").Value = '  -3.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4827'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -1.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9059'
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '108.40'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  -3.88%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.901'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -3.04%  '

$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.680'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -0.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4178'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -0.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.090'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -0.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1246'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +0.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.03'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  -0.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05803'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  -2.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8944'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +0.98%  '
